# Auto-generated: apply numeric recalculation updates to columns H-N
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 100.318184
$ws.Range("I9").Value = 71.42104999999999
$ws.Range("J9").Value = 283.33334
$ws.Range("K9").Value = 71.42104999999999
$ws.Range("L9").Value = 283.33334
$ws.Range("M9").Value = 97.57895000000001
$ws.Range("N9").Value = -621.33334

$ws.Range("H62").Value = 97228460
$ws.Range("I62").Value = 41674744
$ws.Range("J62").Value = 208335920
$ws.Range("K62").Value = 41674744
$ws.Range("L62").Value = 208335920
$ws.Range("M62").Value = -41674120
$ws.Range("N62").Value = -208337168

$ws.Range("H65").Value = 97228460
$ws.Range("I65").Value = 41674744
$ws.Range("J65").Value = 208335920
$ws.Range("K65").Value = 208373720
$ws.Range("L65").Value = 1041679600
$ws.Range("M65").Value = -208370600
$ws.Range("N65").Value = -1041685840

$ws.Range("H86").Value = 1928797.8
$ws.Range("I86").Value = 2700
$ws.Range("J86").Value = 3469676
$ws.Range("K86").Value = 2700
$ws.Range("L86").Value = 3469676
$ws.Range("M86").Value = -1577
$ws.Range("N86").Value = -3471922

$ws.Range("H89").Value = 1928797.8
$ws.Range("I89").Value = 2700
$ws.Range("J89").Value = 3469676
$ws.Range("K89").Value = 13500
$ws.Range("L89").Value = 17348380
$ws.Range("M89").Value = -7884
$ws.Range("N89").Value = -17359612

$ws.Range("H132").Value = 2725645.2
$ws.Range("I132").Value = 558876.6
$ws.Range("J132").Value = 22226564
$ws.Range("K132").Value = 1676629.8
$ws.Range("L132").Value = 66679692
$ws.Range("M132").Value = -1674099.8
$ws.Range("N132").Value = -66684752

$ws.Range("H137").Value = 13151539
$ws.Range("I137").Value = 3379229.8
$ws.Range("J137").Value = 37256572
$ws.Range("K137").Value = 10137689.4
$ws.Range("L137").Value = 111769716
$ws.Range("M137").Value = -10135139.4
$ws.Range("N137").Value = -111774816

$ws.Range("H138").Value = 3715.8354
$ws.Range("I138").Value = 2777.25
$ws.Range("J138").Value = 4231.137
$ws.Range("K138").Value = 8331.75
$ws.Range("L138").Value = 12693.411
$ws.Range("M138").Value = -3191.75
$ws.Range("N138").Value = -22973.411

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 24567.75
$ws.Range("J55").Value = 24567.75
$ws.Range("L55").Value = 24567.75
$ws.Range("N55").Value = -25197.75

$ws.Range("H61").Value = 3243125
$ws.Range("I61").Value = 1603839.8
$ws.Range("J61").Value = 11767408
$ws.Range("K61").Value = 1603839.8
$ws.Range("L61").Value = 11767408
$ws.Range("M61").Value = -1603627.8
$ws.Range("N61").Value = -11767832

$ws.Range("H136").Value = 3243125
$ws.Range("I136").Value = 1603839.8
$ws.Range("J136").Value = 11767408
$ws.Range("K136").Value = 4811519.4
$ws.Range("L136").Value = 35302224
$ws.Range("M136").Value = -4808969.4
$ws.Range("N136").Value = -35307324

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 9500
$ws.Range("J9").Value = 9500
$ws.Range("L9").Value = 9500
$ws.Range("N9").Value = -9836

$ws.Range("H20").Value = 20012774
$ws.Range("I20").Value = 35724332
$ws.Range("J20").Value = 16243.637
$ws.Range("K20").Value = 35724332
$ws.Range("L20").Value = 16243.637
$ws.Range("M20").Value = -35724085
$ws.Range("N20").Value = -16737.637

$ws.Range("H44").Value = 7500
$ws.Range("J44").Value = 7500
$ws.Range("L44").Value = 7500
$ws.Range("N44").Value = -8494

$ws.Range("H134").Value = 12367499
$ws.Range("I134").Value = 13305614
$ws.Range("J134").Value = 172000
$ws.Range("K134").Value = 39916842
$ws.Range("L134").Value = 516000
$ws.Range("M134").Value = -39914307
$ws.Range("N134").Value = -521070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3628055
$ws.Range("I31").Value = 8334884.5
$ws.Range("J31").Value = 7416.769
$ws.Range("K31").Value = 8334884.5
$ws.Range("L31").Value = 7416.769
$ws.Range("M31").Value = -8334589.5
$ws.Range("N31").Value = -8006.769

$ws.Range("H34").Value = 3628055
$ws.Range("I34").Value = 8334884.5
$ws.Range("J34").Value = 7416.769
$ws.Range("K34").Value = 8334884.5
$ws.Range("L34").Value = 7416.769
$ws.Range("M34").Value = -8334682.5
$ws.Range("N34").Value = -7820.769

$ws.Range("H51").Value = 28030.1
$ws.Range("J51").Value = 10037.625
$ws.Range("L51").Value = 10037.625
$ws.Range("N51").Value = -11509.625

$ws.Range("H58").Value = 2279005.2
$ws.Range("I58").Value = 8582.538
$ws.Range("J58").Value = 6495505
$ws.Range("K58").Value = 8582.538
$ws.Range("L58").Value = 6495505
$ws.Range("M58").Value = -8379.538
$ws.Range("N58").Value = -6495911

$ws.Range("H61").Value = 28030.1
$ws.Range("J61").Value = 10037.625
$ws.Range("L61").Value = 10037.625
$ws.Range("N61").Value = -10733.625

$ws.Range("H136").Value = 2279005.2
$ws.Range("I136").Value = 8582.538
$ws.Range("J136").Value = 6495505
$ws.Range("K136").Value = 25747.614
$ws.Range("L136").Value = 19486515
$ws.Range("M136").Value = -23197.614
$ws.Range("N136").Value = -19491615

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 8598.538
$ws.Range("I68").Value = 630
$ws.Range("J68").Value = 10989.1
$ws.Range("K68").Value = 1890
$ws.Range("L68").Value = 32967.3
$ws.Range("M68").Value = -1079
$ws.Range("N68").Value = -34589.3

$ws.Range("H71").Value = 8598.538
$ws.Range("I71").Value = 630
$ws.Range("J71").Value = 10989.1
$ws.Range("K71").Value = 5670
$ws.Range("L71").Value = 98901.90000000001
$ws.Range("M71").Value = -1614
$ws.Range("N71").Value = -107013.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4194977.5
$ws.Range("I70").Value = 1693305.5
$ws.Range("J70").Value = 11908466
$ws.Range("K70").Value = 1693305.5
$ws.Range("L70").Value = 11908466
$ws.Range("M70").Value = -1693035.5
$ws.Range("N70").Value = -11909006

$ws.Range("H73").Value = 4194977.5
$ws.Range("I73").Value = 1693305.5
$ws.Range("J73").Value = 11908466
$ws.Range("K73").Value = 1693305.5
$ws.Range("L73").Value = 11908466
$ws.Range("M73").Value = -1692369.5
$ws.Range("N73").Value = -11910338

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4087561.8
$ws.Range("I132").Value = 5106866.5
$ws.Range("J132").Value = 10342.714
$ws.Range("K132").Value = 15320599.5
$ws.Range("L132").Value = 31028.142
$ws.Range("M132").Value = -15318069.5
$ws.Range("N132").Value = -36088.142

$ws.Range("H136").Value = 5685536.5
$ws.Range("I136").Value = 7356306
$ws.Range("J136").Value = 4920
$ws.Range("K136").Value = 22068918
$ws.Range("L136").Value = 14760
$ws.Range("M136").Value = -22066368
$ws.Range("N136").Value = -19860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 989.5161000000001
$ws.Range("I122").Value = 807.6316
$ws.Range("J122").Value = 1277.5
$ws.Range("K122").Value = 2422.8948
$ws.Range("L122").Value = 3832.5
$ws.Range("M122").Value = 27.10519999999997
$ws.Range("N122").Value = -8732.5

$ws.Range("H132").Value = 663424.5
$ws.Range("I132").Value = 1774
$ws.Range("K132").Value = 5322
$ws.Range("M132").Value = -2792
